# Slide 12 ("Where do I start?(Things I learned)") - Content Placeholder 2:
#   "Smaller code = smaller tests"          -> "Simpler code = simpler tests"
#   "Start with what " + "you know" (2 runs) -> "Start with what you know" (1 run)
#
# NOTE: the diff also re-caches the deck's "update automatically" date field
# (datetimeFigureOut, 3/26/2018 -> 4/2/2018) across the slide master/layouts.
# That text is a cached render of an auto-updating field tied to the
# save-time clock, not an authored edit reachable from the object model, so
# it is intentionally left alone here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Walk paragraphs back-to-front so earlier Start/Length values stay valid
# even if a replacement text has a different length than the original.
$para2 = $tr.Paragraphs(2)
$tr.Characters($para2.Start, $para2.Length).Text = "Start with what you know"

$para1 = $tr.Paragraphs(1)
$tr.Characters($para1.Start, $para1.Length).Text = "Simpler code = simpler tests"
